# Fix the maxMin test file: the O/P columns for rows 2, 4 and 6 on the
# "Template" sheet had their Min/Max values swapped. Restore the correct
# values and leave the active cell selection on P6 (matching what Excel
# would record after the last edit was made there).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Swap O/P values for the affected rows.
foreach ($row in 2, 4, 6) {
    $oCell = $ws.Cells.Item($row, 15)  # column O
    $pCell = $ws.Cells.Item($row, 16)  # column P

    $oVal = $oCell.Value2
    $pVal = $pCell.Value2

    $oCell.Value = $pVal
    $pCell.Value = $oVal
}

# Make the sheet active and select P6, matching the saved cursor position.
$ws.Activate()
$ws.Range("P6").Select()
